# Update "想去人数" (want-to-go count) figures in column F across sheets,
# reflecting refreshed data output generated at 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibitions) ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value  = 3141
$wsExhibit.Range("F8").Value  = 921
$wsExhibit.Range("F16").Value = 7702
$wsExhibit.Range("F17").Value = 336
$wsExhibit.Range("F18").Value = 2452
$wsExhibit.Range("F22").Value = 460
$wsExhibit.Range("F25").Value = 1129
$wsExhibit.Range("F34").Value = 162
$wsExhibit.Range("F35").Value = 271
$wsExhibit.Range("F38").Value = 331
$wsExhibit.Range("F40").Value = 212

# --- Sheet: 演出 (Performances) ---
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F5").Value = 12

# --- Sheet: 全部类型 (All types) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F5").Value  = 3141
$wsAll.Range("F10").Value = 921
$wsAll.Range("F19").Value = 7702
$wsAll.Range("F20").Value = 336
$wsAll.Range("F21").Value = 2452
$wsAll.Range("F22").Value = 12
$wsAll.Range("F26").Value = 460
$wsAll.Range("F29").Value = 1129
$wsAll.Range("F38").Value = 162
$wsAll.Range("F39").Value = 271
$wsAll.Range("F42").Value = 331
$wsAll.Range("F47").Value = 212
